# Edit: update CasesTab query (No -> Yes) and replace the shared StatQuery
# text in the three tab rows (CasesTab/SamplesTab/FilesTab) with the new
# combined statistics query. Also nudge the saved view state (zoom/selection)
# to match the author's final position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Decode-Base64Text([string]$b64) {
    $bytes = [Convert]::FromBase64String($b64)
    return [System.Text.Encoding]::UTF8.GetString($bytes)
}

$casesQueryYes = Decode-Base64Text "TUFUQ0ggKHM6c3R1ZHkpPC1bKl0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYykKCk1BVENIIChjKTwtLShkaWFnOmRpYWdub3NpcykKIE1BVENIIChwOnByb2dyYW0pPC1bKl0tKHM6c3R1ZHkpPC1bKl0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYyksIChjKTwtLShkaWFnOmRpYWdub3NpcykKCVdIRVJFIHMuY2xpbmljYWxfc3R1ZHlfZGVzaWduYXRpb24gSU4gWydVQkMwMSddIGFuZCBkZW1vLm5ldXRlcmVkX2luZGljYXRvciBpbiBbICdZZXMnXSBPUFRJT05BTCBNQVRDSCAoc2FtcDpzYW1wbGUpLS0+KGMpCk9QVElPTkFMIE1BVENIIChjbzpjb2hvcnQpPC1bKl0tKGMpCldJVEggRElTVElOQ1QgYywgcywgZGVtbywgZGlhZywgY28KUkVUVVJOICBjb2FsZXNjZShjLmNhc2VfaWQsICcnKSBBUyBgQ2FzZSBJRGAgLAogICAgICAgIGNvYWxlc2NlKHMuY2xpbmljYWxfc3R1ZHlfZGVzaWduYXRpb24sICcnKSBBUyBgU3R1ZHkgQ29kZWAgLAogICAgICAgIGNvYWxlc2NlKHMuY2xpbmljYWxfc3R1ZHlfdHlwZSwgJycpIEFTICBgU3R1ZHkgVHlwZWAsCiAgICAgICAgY29hbGVzY2UoZGVtby5icmVlZCwgJycpIEFTIEJyZWVkICwKICAgICAgICBjb2FsZXNjZShkaWFnLmRpc2Vhc2VfdGVybSwgJycpIEFTIERpYWdub3NpcyAsCiAgICAgICAgY29hbGVzY2UoZGlhZy5zdGFnZV9vZl9kaXNlYXNlLCAnJykgQVMgYFN0YWdlIG9mIERpc2Vhc2VgICwKICAgICAgICBjb2FsZXNjZShkZW1vLnBhdGllbnRfYWdlX2F0X2Vucm9sbG1lbnQsICcnKSBBUyBBZ2UgLAogICAgICAgIGNvYWxlc2NlKGRlbW8uc2V4LCAnJykgQVMgU2V4ICwKICAgICAgICBjb2FsZXNjZShkZW1vLm5ldXRlcmVkX2luZGljYXRvciwgJycpIEFTIGBOZXV0ZXJlZCBTdGF0dXNgLAogICAgICAgIGNvYWxlc2NlKGRlbW8ud2VpZ2h0LCAnJykgQVMgYFdlaWdodCAoa2cpYCwKICAgICAgICBjb2FsZXNjZShkaWFnLmJlc3RfcmVzcG9uc2UsICcnKSBBUyBgUmVzcG9uc2UgdG8gVHJlYXRtZW50YA=="
$newStatQuery  = Decode-Base64Text "IE1BVENIIChwOnByb2dyYW0pPC0tKHM6c3R1ZHkpPC0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYyksIChjKTwtLShkaWFnOmRpYWdub3NpcykKICAgICAgV0hFUkUgKHNpemUoW10pID0gMCBPUiBzLmNsaW5pY2FsX3N0dWR5X2Rlc2lnbmF0aW9uIElOIFtdKQogICAgICAgIEFORCAocy5zdHVkeV9kaXNwb3NpdGlvbiA9ICdVbnJlc3RyaWN0ZWQnKQogICAgICAgIEFORCAoc2l6ZShbXSkgPSAwIE9SIHMuY2xpbmljYWxfc3R1ZHlfdHlwZSBJTiBbXSkKICAgICAgICBBTkQgKHNpemUoWydVQkMwMSddKSA9IDAgT1IgZGVtby5icmVlZCBJTiBbJ1VCQzAxJ10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZGVtby5zZXggSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZGVtby5uZXV0ZXJlZF9pbmRpY2F0b3IgSU4gW1llc10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZGlhZy5kaXNlYXNlX3Rlcm0gSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZGlhZy5wcmltYXJ5X2Rpc2Vhc2Vfc2l0ZSBJTiBbXSkKICAgICAgICBBTkQgKHNpemUoW10pID0gMCBPUiBkaWFnLnN0YWdlX29mX2Rpc2Vhc2UgSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZGlhZy5iZXN0X3Jlc3BvbnNlIElOIFtdKQogICAgT1BUSU9OQUwgTUFUQ0ggKGMpLS0+KGNvOmNvaG9ydCkKICAgIE9QVElPTkFMIE1BVENIIChmOmZpbGUpLVsqXS0+KGMpCiAgICBPUFRJT05BTCBNQVRDSCAoZiktLT4ocGFyZW50KQogICAgT1BUSU9OQUwgTUFUQ0ggKHNhbXA6c2FtcGxlKS0tPihjKQogICAgT1BUSU9OQUwgTUFUQ0ggKHNhbXApPC0tKGFsOmFsaXF1b3QpCiAgICBXSVRIIERJU1RJTkNUIGMgQVMgYywgcCwgcywgY28sIGRlbW8sIGRpYWcsIGYsIHBhcmVudCwgc2FtcCwgYWwKICAgICAgV0hFUkUgKHNpemUoW10pID0gMCBPUiBzYW1wLnN1bW1hcml6ZWRfc2FtcGxlX3R5cGUgSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1Igc2FtcC5zcGVjaWZpY19zYW1wbGVfcGF0aG9sb2d5IElOIFtdKQogICAgICAgIEFORCAoc2l6ZShbXSkgPSAwIE9SIHNhbXAuc2FtcGxlX3NpdGUgSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgaGVhZChsYWJlbHMocGFyZW50KSkgSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZi5maWxlX3R5cGUgSU4gW10pCiAgICAgICAgQU5EIChzaXplKFtdKSA9IDAgT1IgZi5maWxlX2Zvcm1hdCBJTiBbXSkKICAgIFdJVEggYy5jYXNlX2lkIEFTIGNhc2VfaWQsCiAgICAgICAgIHMuY2xpbmljYWxfc3R1ZHlfZGVzaWduYXRpb24gQVMgc3R1ZHlfY29kZSwKICAgICAgICAgcy5jbGluaWNhbF9zdHVkeV90eXBlIEFTIHN0dWR5X3R5cGUsCiAgICAgICAgIGNvLmNvaG9ydF9kZXNjcmlwdGlvbiBBUyBjb2hvcnQsCiAgICAgICAgIGRlbW8uYnJlZWQgQVMgYnJlZWQsCiAgICAgICAgIGRpYWcuZGlzZWFzZV90ZXJtIEFTIGRpYWdub3NpcywKICAgICAgICAgZGlhZy5zdGFnZV9vZl9kaXNlYXNlIEFTIHN0YWdlX29mX2Rpc2Vhc2UsCiAgICAgICAgIGRpYWcucHJpbWFyeV9kaXNlYXNlX3NpdGUgQVMgZGlzZWFzZV9zaXRlLAogICAgICAgICBkZW1vLnBhdGllbnRfYWdlX2F0X2Vucm9sbG1lbnQgQVMgYWdlLAogICAgICAgICBkZW1vLnNleCBBUyBzZXgsCiAgICAgICAgIGRlbW8ubmV1dGVyZWRfaW5kaWNhdG9yIEFTIG5ldXRlcmVkX3N0YXR1cywKICAgICAgICAgZGVtby53ZWlnaHQgQVMgd2VpZ2h0LAogICAgICAgICBkaWFnLmJlc3RfcmVzcG9uc2UgQVMgcmVzcG9uc2VfdG9fdHJlYXRtZW50LAogICAgICAgICBzYW1wLnNhbXBsZV9pZCBBUyBzYW1wbGVfaWQsCiAgICAgICAgIGYudXVpZCBBUyBmaWxlX2lkLAogICAgICAgICBhbAogICAgUkVUVVJOCkNPVU5UKERJU1RJTkNUIGZpbGVfaWQpIGFzIG51bWJlcl9vZl9maWxlcywKQ09VTlQoRElTVElOQ1Qgc2FtcGxlX2lkKSBhcyBudW1iZXJfb2Zfc2FtcGxlLApDT1VOVChESVNUSU5DVCBjYXNlX2lkKSBhcyBudW1iZXJfb2ZfY2FzZXMsCkNPVU5UKERJU1RJTkNUIHN0dWR5X2NvZGUpIGFzIG51bWJlcl9vZl9zdHVkeSwKQ09VTlQoRElTVElOQ1QgYWwpIGFzIG51bWJlcl9vZl9hbGlxdW90CiAgICA="

# CasesTab row (row 2): query text toggles neutered_indicator filter to 'Yes'
$ws.Range("B2").Value = $casesQueryYes

# StatQuery column (C) for all three tabs now shares the new combined query
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Row heights grow to fit the much longer StatQuery text
$ws.Rows(2).RowHeight = 409.6
$ws.Rows(3).RowHeight = 409.6
$ws.Rows(4).RowHeight = 409.6

# Restore the saved view/selection state
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("B4").Select()
